# fix: improve error handling in file upload and update test data
#
# The "Recipients" sheet's second data row (row 3, Jane Smith) had its
# email address changed from jane.smith@example.com to
# insomewhererandom@gmail.com, and that cell is now turned into a live
# mailto: hyperlink (Excel's default blue/underlined "Hyperlink" style).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newEmail = "insomewhererandom@gmail.com"

# Update the cell text...
$ws.Range("B3").Value = $newEmail

# ...and wire it up as a real mailto: hyperlink (this also applies the
# built-in "Hyperlink" cell style: blue, underlined Calibri).
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:" + $newEmail)

# Match the saved selection/active-cell state on the sheet.
$ws.Range("B3").Select() | Out-Null
